# 996: Add T2A to xlsx files and extract process
$wb = $excel.ActiveWorkbook

# --- Update WMT_Extract sheet view/selection (no longer scrolled/selecting rows 2-3) ---
$wmt = $wb.Worksheets.Item("WMT_Extract")
[void]$wmt.Range("A1:XFD1").Select()

# --- Add the new T2A sheet after the last existing sheet (GS) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$t2a = $wb.Worksheets.Add($null, $lastSheet)
$t2a.Name = "T2A"

# Copy the header row formatting from WMT_Extract so the new sheet matches
# styles (fonts etc.) used on the original extract sheet.
$srcHeader = $wmt.Range("A1:AO1")
$srcHeader.Copy()
$destHeader = $t2a.Range("A1:AO1")
[void]$destHeader.PasteSpecial(-4122)

# Copy the actual header values/shared-strings cell by cell.
for ($col = 1; $col -le 41; $col++) {
  $val = $wmt.Cells.Item(1, $col).Value2
  $t2a.Cells.Item(1, $col).Value2 = $val
}

# Match the source sheet's (taller) header row height.
$t2a.Rows.Item(1).RowHeight = $wmt.Rows.Item(1).RowHeight

[void]$t2a.Range("A1:XFD1").Select()
[void]$t2a.Activate()
